$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Job") values from 4 to 6 for rows 2 through 19
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 3).Value = 6
}

# Update the active selection to E19
$ws.Range("E19").Select()
